$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing note about "_custom" files to also mention standalone creation.
$ws.Range("A4").Value = "Dateien mit Endung ""_custom"" wurden auf Basis anderer Dateien oder eigenständig erstellt"

# Insert a new row above the current row 5 (the empty, bordered divider row),
# shifting it and the header row below it down by one.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row with a note about the "letters" icons.
$ws.Range("A5").Value = "Alle svg Icons im ordner ""letters"" wurden eigenständig erstellt"

# Update the active cell selection as recorded at save time.
$ws.Range("J29").Select()
